$d = $word.ActiveDocument
$d.Paragraphs.Item(1).Range.Find.Execute("2023-04-19 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-04-20 Thursday", 2) | Out-Null
$t = $d.Tables.Item(1)
$t.Cell(1,1).Range.Text = "52×52="
$t.Cell(1,2).Range.Text = "21×22="
$t.Cell(1,3).Range.Text = "90×51="
$t.Cell(1,4).Range.Text = "55×80="
$t.Cell(1,5).Range.Text = "19×54="
$t.Cell(2,1).Range.Text = "40×41="
$t.Cell(2,2).Range.Text = "99×80="
$t.Cell(2,3).Range.Text = "88×75="
$t.Cell(2,4).Range.Text = "35×13="
$t.Cell(2,5).Range.Text = "13×14="
$t.Cell(3,1).Range.Text = "21×74="
$t.Cell(3,2).Range.Text = "87×82="
$t.Cell(3,3).Range.Text = "27×27="
$t.Cell(3,4).Range.Text = "18×29="
$t.Cell(3,5).Range.Text = "99×88="
$t.Cell(4,1).Range.Text = "52×100="
$t.Cell(4,2).Range.Text = "69×21="
$t.Cell(4,3).Range.Text = "32×81="
$t.Cell(4,4).Range.Text = "33×18="
$t.Cell(4,5).Range.Text = "70×54="
$t.Cell(5,1).Range.Text = "48×51="
$t.Cell(5,2).Range.Text = "12×67="
$t.Cell(5,3).Range.Text = "25×85="
$t.Cell(5,4).Range.Text = "79×48="
$t.Cell(5,5).Range.Text = "81×23="
$t.Cell(6,1).Range.Text = "32×74="
$t.Cell(6,2).Range.Text = "45×26="
$t.Cell(6,3).Range.Text = "86×55="
$t.Cell(6,4).Range.Text = "34×35="
$t.Cell(6,5).Range.Text = "42×14="
$t.Cell(7,1).Range.Text = "92×51="
$t.Cell(7,2).Range.Text = "59×33="
$t.Cell(7,3).Range.Text = "93×85="
$t.Cell(7,4).Range.Text = "74×25="
$t.Cell(7,5).Range.Text = "56×65="
$t.Cell(8,1).Range.Text = "20×81="
$t.Cell(8,2).Range.Text = "13×82="
$t.Cell(8,3).Range.Text = "11×92="
$t.Cell(8,4).Range.Text = "19×50="
$t.Cell(8,5).Range.Text = "91×97="
$t.Cell(9,1).Range.Text = "72×74="
$t.Cell(9,2).Range.Text = "39×50="
$t.Cell(9,3).Range.Text = "79×18="
$t.Cell(9,4).Range.Text = "84×13="
$t.Cell(9,5).Range.Text = "52×47="
$t.Cell(10,1).Range.Text = "18×86="
$t.Cell(10,2).Range.Text = "29×64="
$t.Cell(10,3).Range.Text = "21×31="
$t.Cell(10,4).Range.Text = "72×45="
$t.Cell(10,5).Range.Text = "88×97="
$t.Cell(11,1).Range.Text = "74×78="
$t.Cell(11,2).Range.Text = "88×14="
$t.Cell(11,3).Range.Text = "87×77="
$t.Cell(11,4).Range.Text = "56×16="
$t.Cell(11,5).Range.Text = "71×81="
$t.Cell(12,1).Range.Text = "91×25="
$t.Cell(12,2).Range.Text = "79×61="
$t.Cell(12,3).Range.Text = "94×27="
$t.Cell(12,4).Range.Text = "36×62="
$t.Cell(12,5).Range.Text = "88×70="
$t.Cell(13,1).Range.Text = "26×88="
$t.Cell(13,2).Range.Text = "91×92="
$t.Cell(13,3).Range.Text = "22×38="
$t.Cell(13,4).Range.Text = "82×56="
$t.Cell(13,5).Range.Text = "66×11="
$t.Cell(14,1).Range.Text = "65×76="
$t.Cell(14,2).Range.Text = "82×17="
$t.Cell(14,3).Range.Text = "13×88="
$t.Cell(14,4).Range.Text = "42×14="
$t.Cell(14,5).Range.Text = "12×94="
$t.Cell(15,1).Range.Text = "81×90="
$t.Cell(15,2).Range.Text = "55×34="
$t.Cell(15,3).Range.Text = "23×58="
$t.Cell(15,4).Range.Text = "60×58="
$t.Cell(15,5).Range.Text = "16×61="
$t.Cell(16,1).Range.Text = "87×58="
$t.Cell(16,2).Range.Text = "65×91="
$t.Cell(16,3).Range.Text = "51×58="
$t.Cell(16,4).Range.Text = "67×95="
$t.Cell(16,5).Range.Text = "84×34="
$t.Cell(17,1).Range.Text = "64×100="
$t.Cell(17,2).Range.Text = "26×26="
$t.Cell(17,3).Range.Text = "54×40="
$t.Cell(17,4).Range.Text = "95×25="
$t.Cell(17,5).Range.Text = "82×29="
$t.Cell(18,1).Range.Text = "63×48="
$t.Cell(18,2).Range.Text = "85×30="
$t.Cell(18,3).Range.Text = "24×71="
$t.Cell(18,4).Range.Text = "78×86="
$t.Cell(18,5).Range.Text = "23×21="
$t.Cell(19,1).Range.Text = "14×39="
$t.Cell(19,2).Range.Text = "67×79="
$t.Cell(19,3).Range.Text = "39×18="
$t.Cell(19,4).Range.Text = "32×23="
$t.Cell(19,5).Range.Text = "93×39="
$t.Cell(20,1).Range.Text = "37×48="
$t.Cell(20,2).Range.Text = "70×29="
$t.Cell(20,3).Range.Text = "41×86="
$t.Cell(20,4).Range.Text = "59×88="
$t.Cell(20,5).Range.Text = "19×27="
Write-Output "done"